$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.214.46"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3
$ws.Range("D3").Value = "2.478.69"

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6
$ws.Range("E6").Value = "  +0.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("D9").Value = "2.478.68"
$ws.Range("E9").Value = "  +1.69%  "

# Row 10
$ws.Range("E10").Value = "  +0.45%  "

# Row 11
$ws.Range("E11").Value = "  +1.73%  "

# Row 12
$ws.Range("E12").Value = "  +0.58%  "

# Row 15
$ws.Range("E15").Value = "  +1.15%  "

# Row 16
$ws.Range("D16").Value = "2.928.61"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17
$ws.Range("D17").Value = "63.127.44"

# Row 18
$ws.Range("D18").Value = "2.477.70"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.45%  "

# Row 20
$ws.Range("E20").Value = "  +1.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.78%  "

# Row 22
$ws.Range("E22").Value = "  +9.01%  "

# Row 24
$ws.Range("E24").Value = "  +0.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$ws.Range("E26").Value = "  +15.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "657.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.66%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0991"
$ws.Range("E28").Value = "  +1.29%  "

# Row 29
$ws.Range("D29").Value = "2.608.50"
$ws.Range("E29").Value = "  +1.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +272.22%  "

# Row 31
$ws.Range("E31").Value = "  +4.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "

# Row 33
$ws.Range("E33").Value = "  +0.98%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.38%  "

# Row 35
$ws.Range("E35").Value = "  +3.91%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("E37").Value = "  +0.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.39%  "

# Row 39
$ws.Range("E39").Value = "  -0.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.95%  "

# Row 43
$ws.Range("E43").Value = "  +0.63%  "

# Row 44
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0320"
$ws.Range("E44").Value = "  -47.37%  "

# Row 45
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.62%  "

# Row 47
$ws.Range("E47").Value = "  +3.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "

# Row 50
$ws.Range("E50").Value = "  +1.95%  "

# Row 51
$ws.Range("E51").Value = "  +0.30%  "
